$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column B and C
$ws.Range("B1").Value = 0.875
$ws.Range("C1").Value = 0.04146758932281736
$ws.Range("B2").Value = 0.9
$ws.Range("C2").Value = 0.01989412373158344

# Add new values in columns D and E
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 202.0454051677636
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 183.0146123922239
